# Insert a new data row at row 692 (pushing the existing 692:733 block down
# to 693:734) and populate it with the new observation:
#   2026/01/22  木  13  201
#
# This mirrors the diff: dimension grows from A1:D733 to A1:D734, a new
# row 692 is inserted, and everything that used to live at rows 692-733
# now lives at rows 693-734 (values unchanged, just shifted down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 692:733 down to 693:734, opening up a blank row 692.
$ws.Rows("692:692").Insert()

# Force text storage for the date/weekday columns so Excel doesn't
# reinterpret "2026/01/22" as a date serial - the rest of the column
# is stored as plain text too.
$ws.Range("A692:B692").NumberFormat = "@"

$ws.Range("A692").Value = "2026/01/22"
$ws.Range("B692").Value = "木"
$ws.Range("C692").Value = 13
$ws.Range("D692").Value = 201

# Drop the temporary text formatting so the new cells end up styleless,
# matching the rest of the data rows in the sheet.
$ws.Range("A692:D692").ClearFormats()
